$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1888
$ws1.Range("F6").Value = 13339
$ws1.Range("F7").Value = 13204
$ws1.Range("F11").Value = 559
$ws1.Range("F13").Value = 680
$ws1.Range("F20").Value = 261
$ws1.Range("F23").Value = 760
$ws1.Range("F24").Value = 18

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 124
$ws2.Range("F9").Value = 32

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 195
$ws3.Range("F3").Value = 44

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 195
$ws4.Range("F3").Value = 1888
$ws4.Range("F8").Value = 13339
$ws4.Range("F9").Value = 13204
$ws4.Range("F13").Value = 559
$ws4.Range("F15").Value = 680
$ws4.Range("F25").Value = 44
$ws4.Range("F27").Value = 261
$ws4.Range("F30").Value = 760
$ws4.Range("F31").Value = 124
$ws4.Range("F33").Value = 18
$ws4.Range("F34").Value = 32
